$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.187.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.369.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.86%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.370.43'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("E11").Value = '  +3.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.950.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.369.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.248.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.52%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '378.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.65%  '
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  +1.64%  '
$ws.Range("E27").Value = '  +11.41%  '
$ws.Range("E28").Value = '  +14.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.21'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.14%  '
$ws.Range("E32").Value = '  +4.28%  '
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.400.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.22%  '
$ws.Range("E37").Value = '  +6.06%  '
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("E39").Value = '  +4.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.05%  '
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.81%  '
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("E47").Value = '  +6.97%  '
$ws.Range("E48").Value = '  +2.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.57%  '
$ws.Range("E50").Value = '  +7.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.322.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.18%  '
